$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply the new font formatting (Calibri 10, theme text color) to row 2 (A2:H2) ---
$rowRange = $ws.Range("A2:H2")
$rowRange.Font.Name = "Calibri"
$rowRange.Font.Size = 10
$rowRange.Font.ThemeColor = 1

# B2 has no content/formatting in the source row -- remove it entirely so it
# does not appear as a styled-but-empty cell.
$ws.Range("B2").Clear()

# --- Write the new data row values (row 2) ---
$ws.Range("A2").Value = "MCH138-1"
$ws.Range("C2").Value = "BOOKS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 21M | GRAP COUNT NUMER: NONE"

# D2 and H2 stay blank (formatted only, no content), matching the source row.

# --- Restore freeze panes (top row frozen) and selection on row 2 ---
$ws.Range("A2:K2").Select()
$excel.ActiveWindow.FreezePanes = $true
